$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.158.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.571.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.575"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.967.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("E14").Value = "  -4.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.579.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.846"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.140.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("E19").Value = "  +4.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.85%  "
$ws.Range("D21").Value = "0.0₃0963"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0810"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("E39").Value = "  +5.62%  "
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.31%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.995.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.818.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.68%  "
